$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIS - Pasep e Cofins")
$ws.Activate()

# New coronavirus transmission / incubation / treatment Q&A pairs
$ws.Range("B7").Value = 'Como o coronavírus é transmitido?'
$ws.Range("C7").Value = 'As investigações sobre as formas de transmissão do coronavírus ainda estão em andamento, mas a disseminação de pessoa para pessoa, ou seja, a contaminação por gotículas respiratórias ou contato, está ocorrendo. <br /> Qualquer pessoa que tenha contato próximo (cerca de 1m) com alguém com sintomas respiratórios está em risco de ser exposta à infecção. <br /> É importante observar que a disseminação de pessoa para pessoa pode ocorrer de forma continuada. <br /> Alguns vírus são altamente contagiosos (como sarampo), enquanto outros são menos. Ainda não está claro com que facilidade o coronavírus se espalha de pessoa para pessoa.'
$ws.Range("B8").Value = 'Quais as formas de transmissão? Como  transmite Como pegar Transmissão'
$ws.Range("C8").Value = ' gotículas de saliva; <br /> <br />espirro; <br /> <br />tosse; <br /> <br />catarro; <br /> <br />contato pessoal próximo, como toque ou aperto de mão; <br /> <br />contato com objetos ou superfícies contaminadas, seguido de contato com a boca, nariz ou olhos.'
$ws.Range("B9").Value = 'Qual o período de incubação ? Tempo de recuperação incubação Quanto demora para recuperar Quantos dias fico em casa Quantos dias para curar'
$ws.Range("C9").Value = 'O período médio de incubação por coronavírus é de 5 dias, com intervalos que chegam a 12 dias, período em que os primeiros sintomas levam para aparecer desde a infecção.'
$ws.Range("B10").Value = 'Qual o tratamento? Tratamento remédio como curar tratamento cura tempo'
$ws.Range("C10").Value = 'Não existe tratamento específico para infecções causadas por coronavírus humano. No caso do coronavírus é indicado repouso e consumo de bastante água, além de algumas medidas adotadas para aliviar os sintomas, conforme cada caso, como, por exemplo: <br /> Uso de medicamento para dor e febre (antitérmicos e analgésicos). <br /> Uso de umidificador no quarto ou tomar banho quente para auxiliar no alívio da dor de garanta e tosse. <br /> Assim que os primeiros sintomas surgirem, é fundamental procurar ajuda médica imediata para confirmar diagnóstico e iniciar o tratamento. <br /> Todos os pacientes que receberem alta durante os primeiros 07 dias do início do quadro (qualquer sintoma independente de febre), devem ser alertados para a possibilidade de piora tardia do quadro clínico e sinais de alerta de complicações como: aparecimento de febre (podendo haver casos iniciais sem febre), elevação ou reaparecimento de febre ou sinais respiratórios, taquicardia (aumento dos batimentos cardíacos), dor pleurítica (dor no peito), fadiga (cansaço) e dispnéia (falta de ar).'

# Reset the remembered selection to B19 (bot returns to initial state)
$ws.Range("B19").Select()
